$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing row 74 values (MV data refresh)
$ws.Cells.Item(74, 2).Value = 24767
$ws.Cells.Item(74, 7).Value = 16536
$ws.Cells.Item(74, 8).Value = 10947

# Add new row 75 with the new quarter's data.
# Force column A to be treated as text so "01-04-2021" is not
# auto-converted to a date serial number, then clear the formatting
# so no extra style gets attached to the cell.
$ws.Cells.Item(75, 1).NumberFormat = "@"
$ws.Cells.Item(75, 1).Value = "01-04-2021"
$ws.Cells.Item(75, 1).ClearFormats()

$ws.Cells.Item(75, 2).Value = 24276
$ws.Cells.Item(75, 3).Value = 8446
$ws.Cells.Item(75, 4).Value = 1044
$ws.Cells.Item(75, 5).Value = 4797
$ws.Cells.Item(75, 6).Value = 2605
$ws.Cells.Item(75, 7).Value = 15830
$ws.Cells.Item(75, 8).Value = 11284
$ws.Cells.Item(75, 9).Value = 4546
